# Updates cryptos list price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.935.86"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.783.24"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.29"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.544"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.07"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("E9").Value = "  +3.01%  "
$ws.Range("E10").Value = "  -4.28%  "
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "2.039.77"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.14"
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("D14").Value = "1.788.09"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "33.911.76"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.50"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.96"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "0.0₃0771"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.64"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.79"
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "1.394.15"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.646"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("E39").Value = "  +8.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.61"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.75"
$ws.Range("E43").Value = "  +15.14%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("E45").Value = "  +10.48%  "
$ws.Range("E46").Value = "  +3.46%  "
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.91"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.34"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "1.942.05"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  +0.20%  "
